# Natmi following Dr Hou advice
# Rebuild Ntrk3-Ptprs LR-pair rows: full 3x3 sending x target cluster cross
# product (ECs, FAPs, sCs) with recomputed specificity/weight statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntrk3"
$ws.Range("C2").Value = "Ptprs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.138349
$ws.Range("H2").Value = 0.415047
$ws.Range("I2").Value = 0.04449686021752534
$ws.Range("J2").Value = 0.04449686021752534
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.090355666666666
$ws.Range("N2").Value = 9.271066999999999
$ws.Range("O2").Value = 0.06928583878088775
$ws.Range("P2").Value = 0.06928583878088775
$ws.Range("Q2").Value = 0.4275476161276666
$ws.Range("R2").Value = 3.847928545148999
$ws.Range("S2").Value = 0.003083002283287159
$ws.Range("T2").Value = 0.003083002283287158

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntrk3"
$ws.Range("C3").Value = "Ptprs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.138349
$ws.Range("H3").Value = 0.415047
$ws.Range("I3").Value = 0.04449686021752534
$ws.Range("J3").Value = 0.04449686021752534
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 25.17096033333333
$ws.Range("N3").Value = 75.512881
$ws.Range("O3").Value = 0.5643334579338453
$ws.Range("P3").Value = 0.5643334579338454
$ws.Range("Q3").Value = 3.482377191156333
$ws.Range("R3").Value = 31.341394720407
$ws.Range("S3").Value = 0.02511106699375503
$ws.Range("T3").Value = 0.02511106699375503

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntrk3"
$ws.Range("C4").Value = "Ptprs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.138349
$ws.Range("H4").Value = 0.415047
$ws.Range("I4").Value = 0.04449686021752534
$ws.Range("J4").Value = 0.04449686021752534
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 16.34167533333333
$ws.Range("N4").Value = 49.025026
$ws.Range("O4").Value = 0.366380703285267
$ws.Range("P4").Value = 0.366380703285267
$ws.Range("Q4").Value = 2.260854440691333
$ws.Range("R4").Value = 20.347689966222
$ws.Range("S4").Value = 0.01630279094048315
$ws.Range("T4").Value = 0.01630279094048315

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntrk3"
$ws.Range("C5").Value = "Ptprs"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.468673666666667
$ws.Range("H5").Value = 7.406021
$ws.Range("I5").Value = 0.7939936469967431
$ws.Range("J5").Value = 0.793993646996743
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.090355666666666
$ws.Range("N5").Value = 9.271066999999999
$ws.Range("O5").Value = 0.06928583878088775
$ws.Range("P5").Value = 0.06928583878088775
$ws.Range("Q5").Value = 7.62907965493411
$ws.Range("R5").Value = 68.661716894407
$ws.Range("S5").Value = 0.05501251581886544
$ws.Range("T5").Value = 0.05501251581886543

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntrk3"
$ws.Range("C6").Value = "Ptprs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.468673666666667
$ws.Range("H6").Value = 7.406021
$ws.Range("I6").Value = 0.7939936469967431
$ws.Range("J6").Value = 0.793993646996743
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 25.17096033333333
$ws.Range("N6").Value = 75.512881
$ws.Range("O6").Value = 0.5643334579338453
$ws.Range("P6").Value = 0.5643334579338454
$ws.Range("Q6").Value = 62.13888693961121
$ws.Range("R6").Value = 559.2499824565009
$ws.Range("S6").Value = 0.4480771803871769
$ws.Range("T6").Value = 0.448077180387177

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntrk3"
$ws.Range("C7").Value = "Ptprs"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.468673666666667
$ws.Range("H7").Value = 7.406021
$ws.Range("I7").Value = 0.7939936469967431
$ws.Range("J7").Value = 0.793993646996743
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 16.34167533333333
$ws.Range("N7").Value = 49.025026
$ws.Range("O7").Value = 0.366380703285267
$ws.Range("P7").Value = 0.366380703285267
$ws.Range("Q7").Value = 40.34226356461622
$ws.Range("R7").Value = 363.080372081546
$ws.Range("S7").Value = 0.2909039507907007
$ws.Range("T7").Value = 0.2909039507907007

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ntrk3"
$ws.Range("C8").Value = "Ptprs"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5021629999999999
$ws.Range("H8").Value = 1.506489
$ws.Range("I8").Value = 0.1615094927857316
$ws.Range("J8").Value = 0.1615094927857315
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.090355666666666
$ws.Range("N8").Value = 9.271066999999999
$ws.Range("O8").Value = 0.06928583878088775
$ws.Range("P8").Value = 0.06928583878088775
$ws.Range("Q8").Value = 1.551862272640333
$ws.Range("R8").Value = 13.966760453763
$ws.Range("S8").Value = 0.01119032067873515
$ws.Range("T8").Value = 0.01119032067873515

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ntrk3"
$ws.Range("C9").Value = "Ptprs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5021629999999999
$ws.Range("H9").Value = 1.506489
$ws.Range("I9").Value = 0.1615094927857316
$ws.Range("J9").Value = 0.1615094927857315
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 25.17096033333333
$ws.Range("N9").Value = 75.512881
$ws.Range("O9").Value = 0.5643334579338453
$ws.Range("P9").Value = 0.5643334579338454
$ws.Range("Q9").Value = 12.63992495386766
$ws.Range("R9").Value = 113.759324584809
$ws.Range("S9").Value = 0.09114521055291333
$ws.Range("T9").Value = 0.09114521055291333

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ntrk3"
$ws.Range("C10").Value = "Ptprs"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5021629999999999
$ws.Range("H10").Value = 1.506489
$ws.Range("I10").Value = 0.1615094927857316
$ws.Range("J10").Value = 0.1615094927857315
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 16.34167533333333
$ws.Range("N10").Value = 49.025026
$ws.Range("O10").Value = 0.366380703285267
$ws.Range("P10").Value = 0.366380703285267
$ws.Range("Q10").Value = 8.206184710412664
$ws.Range("R10").Value = 73.85566239371398
$ws.Range("S10").Value = 0.05917396155408308
$ws.Range("T10").Value = 0.05917396155408308

